# Generate Report for Handback
# - Status text updated from "Ready for handoff" to "Handed back: in sync with en-US"
#   (Overview!E2:F3, zh-cn!C2:C3, de-de!C2:C3)
# - zh-cn / de-de sheets: populate "Latest Target File" (I) and
#   "Latest Handback File" (J) columns with the handed-back file names, and
#   refresh the "Latest Handback DateTime" (K) column with the handback time.
# - Column widths widened to fit the new, longer content.

$wb = $excel.ActiveWorkbook

$statusText = "Handed back: in sync with en-US"

# Cornflower-blue underlined colour used by the existing "HyperLink" style
# in this workbook (rgb FF6495ED == RGB(100,149,237)).
$hyperlinkColor = 15570276

# ColumnWidth inputs chosen so the persisted <col width> lands on the
# nearest representable grid value to the widths used by the target sheet.
$wideStatusColWidth = 29.16666676666667   # -> stored width ~29.98 (was ~17.22)
$wideFileColWidth   = 39.166666766666665  # -> stored width 40 (was ~18.65 / ~21.71)

# ---------------------------------------------------------------------
# Overview sheet: status column (E & F) text + width
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $statusText
$wsOverview.Range("F2").Value = $statusText
$wsOverview.Range("E3").Value = $statusText
$wsOverview.Range("F3").Value = $statusText
$wsOverview.Columns.Item(5).ColumnWidth = $wideStatusColWidth
$wsOverview.Columns.Item(6).ColumnWidth = $wideStatusColWidth

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("C2").Value = $statusText
$wsZhCn.Range("C3").Value = $statusText
$wsZhCn.Columns.Item(3).ColumnWidth = $wideStatusColWidth
$wsZhCn.Columns.Item(9).ColumnWidth = $wideFileColWidth
$wsZhCn.Columns.Item(10).ColumnWidth = $wideFileColWidth

$md2016Target = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d086054d155fe516891a612d645e1042bf90dda0/e2e/c316ba12-b630-4122-9583-8507e8eee9ef.md"

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I2"), $md2016Target, "", "", "c316ba12-b630-4122-9583-8507e8eee9ef.md")
$wsZhCn.Range("I2").Font.Underline = $true
$wsZhCn.Range("I2").Font.Color = $hyperlinkColor
$wsZhCn.Range("J2").Value = "c316ba12-b630-4122-9583-8507e8eee9ef.fb8b5184faab6789b8c930836fb6e0ce9677c175.zh-cn.xlf"
$wsZhCn.Range("K2").Value = "2016-08-17 03:04:50"

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I3"), $md2016Target, "", "", "c316ba12-b630-4122-9583-8507e8eee9ef.md")
$wsZhCn.Range("I3").Font.Underline = $true
$wsZhCn.Range("I3").Font.Color = $hyperlinkColor
$wsZhCn.Range("J3").Value = "c316ba12-b630-4122-9583-8507e8eee9ef.fb8b5184faab6789b8c930836fb6e0ce9677c175.zh-cn.xlf"
$wsZhCn.Range("K3").Value = "2016-08-17 03:04:50"

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("C2").Value = $statusText
$wsDeDe.Range("C3").Value = $statusText
$wsDeDe.Columns.Item(3).ColumnWidth = $wideStatusColWidth
$wsDeDe.Columns.Item(9).ColumnWidth = $wideFileColWidth
$wsDeDe.Columns.Item(10).ColumnWidth = $wideFileColWidth

$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I2"), $md2016Target, "", "", "c316ba12-b630-4122-9583-8507e8eee9ef.md")
$wsDeDe.Range("I2").Font.Underline = $true
$wsDeDe.Range("I2").Font.Color = $hyperlinkColor
$wsDeDe.Range("J2").Value = "c316ba12-b630-4122-9583-8507e8eee9ef.fb8b5184faab6789b8c930836fb6e0ce9677c175.de-de.xlf"
$wsDeDe.Range("K2").Value = "2016-08-17 03:04:58"

$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I3"), $md2016Target, "", "", "c316ba12-b630-4122-9583-8507e8eee9ef.md")
$wsDeDe.Range("I3").Font.Underline = $true
$wsDeDe.Range("I3").Font.Color = $hyperlinkColor
$wsDeDe.Range("J3").Value = "c316ba12-b630-4122-9583-8507e8eee9ef.fb8b5184faab6789b8c930836fb6e0ce9677c175.de-de.xlf"
$wsDeDe.Range("K3").Value = "2016-08-17 03:04:58"
